$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (new Price (D) value or $null if unchanged, new Volume(1h) (E) value)
$updates = @(
    @{ Row = 2;  D = "60.449.10";  E = "  -3.91%  " },
    @{ Row = 3;  D = "3.308.96";   E = "  -4.15%  " },
    @{ Row = 4;  D = $null;        E = "  -0.11%  " },
    @{ Row = 5;  D = "558.67";     E = "  -3.75%  " },
    @{ Row = 6;  D = "143.37";     E = "  -4.70%  " },
    @{ Row = 7;  D = $null;        E = "  -0.02%  " },
    @{ Row = 8;  D = "3.311.23";   E = "  -4.09%  " },
    @{ Row = 9;  D = $null;        E = "  -2.24%  " },
    @{ Row = 10; D = $null;        E = "  -2.52%  " },
    @{ Row = 11; D = $null;        E = "  -3.96%  " },
    @{ Row = 12; D = $null;        E = "  -1.95%  " },
    @{ Row = 13; D = $null;        E = "  -4.37%  " },
    @{ Row = 14; D = $null;        E = "  +0.79%  " },
    @{ Row = 15; D = "27.23";      E = "  -3.94%  " },
    @{ Row = 16; D = $null;        E = "  -4.45%  " },
    @{ Row = 17; D = $null;        E = "  -3.97%  " },
    @{ Row = 18; D = "60.358.18";  E = "  -4.16%  " },
    @{ Row = 19; D = $null;        E = "  -4.81%  " },
    @{ Row = 20; D = $null;        E = "  -1.58%  " },
    @{ Row = 21; D = $null;        E = "  -4.74%  " },
    @{ Row = 22; D = "375.82";     E = "  -3.23%  " },
    @{ Row = 23; D = "73.95";      E = "  -1.67%  " },
    @{ Row = 24; D = $null;        E = "  -4.18%  " },
    @{ Row = 25; D = $null;        E = "  +0.09%  " },
    @{ Row = 26; D = "3.436.64";   E = "  -4.32%  " },
    @{ Row = 27; D = $null;        E = "  -8.84%  " },
    @{ Row = 28; D = $null;        E = "  -7.32%  " },
    @{ Row = 29; D = $null;        E = "  -0.57%  " },
    @{ Row = 30; D = $null;        E = "  -6.00%  " },
    @{ Row = 31; D = $null;        E = "  -0.12%  " },
    @{ Row = 32; D = $null;        E = "  -4.64%  " },
    @{ Row = 33; D = $null;        E = "  -4.34%  " },
    @{ Row = 34; D = $null;        E = "  -3.28%  " },
    @{ Row = 35; D = $null;        E = "  -6.30%  " },
    @{ Row = 36; D = "5.24";       E = "  -3.95%  " },
    @{ Row = 37; D = "166.56";     E = "  -1.62%  " },
    @{ Row = 38; D = $null;        E = "  -7.30%  " },
    @{ Row = 39; D = "6.74";       E = "  -3.31%  " },
    @{ Row = 40; D = $null;        E = "  -15.71%  " },
    @{ Row = 41; D = $null;        E = "  -4.69%  " },
    @{ Row = 42; D = $null;        E = "  -5.57%  " },
    @{ Row = 43; D = "41.94";      E = "  -1.83%  " },
    @{ Row = 44; D = $null;        E = "  -4.42%  " },
    @{ Row = 45; D = $null;        E = "  -4.60%  " },
    @{ Row = 46; D = $null;        E = "  -6.49%  " },
    @{ Row = 47; D = "1.12";       E = "  -5.03%  " },
    @{ Row = 48; D = "2.357.20";   E = "  -7.75%  " },
    @{ Row = 49; D = $null;        E = "  -0.15%  " },
    @{ Row = 50; D = $null;        E = "  -6.02%  " },
    @{ Row = 51; D = $null;        E = "  -4.05%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        # Some "Price" values look like plain numbers (e.g. "558.67").
        # Force them to be stored as text so they match the original
        # inline-string / text representation instead of being
        # auto-converted to a floating point number by Excel.
        $looksNumeric = $u.D -match '^-?\d+(\.\d+)?$'
        if ($looksNumeric) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    $ws.Cells.Item($r, 5).Value = $u.E
}
